$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$questionsText = @'
questions = [
    {
        "title": "You are running a web application on GCP that is logging to Cloud Logging. The application is configured to log any errors. To ensure that your application is providing a smooth user experience, you want to know if there are more than 10 errors per minute in the logs. Which of the following is the best solution?",
        "ques_type": 2,
        "options": [
            "Write a Cloud Function that will get the logs from Cloud Logging and perform a regex for errors. ",
            "Create a Compute Engine instance to run a custom Python script that will get the logs from Cloud Logging and perform a regex for errors. ",
            "Create a Cloud Logging Logs-based metric with an alerting policy and push it to Pub/Sub subscription.",
            "Monitor the VPC network Flow logging to run a custom Python script that will get the logs from Cloud Logging and perform a regex for errors."
        ],
        "score": "Create a Cloud Logging Logs-based metric with an alerting policy and push it to Pub/Sub subscription."
    },
    {
        "title": "You are in charge of managing a fleet of servers both in GCP and in an on-premise data center. After a security audit, you noticed the servers have not been patched or updated in a while and are now susceptible to common vulnerabilities and exposures (CVEs).Which of the following is the most suitable solution to easily update and patch your servers in a consistent manner with minimum operational overhead?",
        "ques_type": 2,
        "options": [
            "Write a Python script to install any updates or patches. Schedule the script to run via a cron scheduler.",
            "Write a Cloud Function that will Secure Shell (SSH) onto all of your servers and execute any updates or patches. Run the Cloud Function every weekend to minimize service disruptions.",
            "Install the GCP OS patch management agent on your GCP and on-premise servers. Secure Shell (SSH) onto your servers. Execute any updates or patches. ",
            "Install the GCP OS patch management agent on your GCP and on-premise servers and deploy the patches via patch jobs. "
        ],
        "score": "Install the GCP OS patch management agent on your GCP and on-premise servers and deploy the patches via patch jobs."
    },
    {
        "title": "Your company has branches globally and has its own private DNS server for its branches as well as contractors. Recently, your team has decided to move the DNS management from local premises to GCP.Which of the following features best suits this requirement?",
        "ques_type": 2,
        "options": [
            "Private Zones",
            "DNS Forwarding",
            "Cloud Logging",
            "Authoritative DNS lookup"
        ],
        "score": "Private Zones"
    },
    {
        "title": "Your team is storing a lot of data in Cloud Storage for backup purposes. The data is rarely accessed and is not considered highly critical, so it does not have to be highly available. You notice that as the volume of data is increasing at a rapid rate, so does the cost.Which of the following solutions minimizes cost?",
        "ques_type": 2,
        "options": [
            "Cloud Storage One Zone - Coldline Storage",
            "Cloud Storage Standard",
            "Cloud Storage Standard - Coldline Storage",
            "Persistent disk"
        ],
        "score": "Cloud Storage One Zone - Coldline Storage"
    }
]
'@

$ws.Rows(1).Delete()
$ws.Range("A1").Value = $questionsText
$ws.Rows(1).EntireRow.AutoFit()
